# New data: 9 May 2022
# Adds the latest monthly unemployment observations (date serial 44652,
# i.e. 2022-04-01 reference month published 9 May 2022) to both sheets:
#   - "Canada"   (sheet1): one new row (national figure)
#   - "Province" (sheet2): ten new rows (one per province/territory group)

$wb = $excel.ActiveWorkbook
$wsCanada   = $wb.Worksheets.Item("Canada")
$wsProvince = $wb.Worksheets.Item("Province")

$newDate = 44652
$dateFmt = "d-mmm-yy"

# ---------------------------------------------------------------------
# Sheet "Canada": append row 29
# ---------------------------------------------------------------------
$rowCanada = 29

$wsCanada.Cells.Item($rowCanada, 1).Value = $newDate
$wsCanada.Cells.Item($rowCanada, 1).NumberFormat = $dateFmt

$wsCanada.Cells.Item($rowCanada, 2).Value = "Canada"
$wsCanada.Cells.Item($rowCanada, 2).NumberFormat = $dateFmt

$wsCanada.Cells.Item($rowCanada, 4).Value = 1085.8
$wsCanada.Cells.Item($rowCanada, 5).Value = 1166.9000000000001

$wsCanada.Range("C" + $rowCanada).Formula = "=(D" + $rowCanada + "-E" + $rowCanada + ")/E" + $rowCanada + "*100"

# ---------------------------------------------------------------------
# Sheet "Province": append rows 272-281
# ---------------------------------------------------------------------
$provinceRows = @(
    @{ Name = "Newfoundland & Labrador"; D = 27.6;   E = 32.299999999999997; First = $true  },
    @{ Name = "Prince Edward Island";    D = 7.5;    E = 7.5;                First = $false },
    @{ Name = "Nova Scotia";             D = 30.6;   E = 34.6;               First = $false },
    @{ Name = "New Brunswick";           D = 28;     E = 32.299999999999997; First = $false },
    @{ Name = "Quebec";                  D = 179.6;  E = 226;                First = $false },
    @{ Name = "Ontario";                 D = 441.6;  E = 468.8;              First = $false },
    @{ Name = "Manitoba";                D = 35.200000000000003; E = 36.299999999999997; First = $false },
    @{ Name = "Saskatchewan";            D = 33.299999999999997; E = 34.200000000000003; First = $false },
    @{ Name = "Alberta";                 D = 146.80000000000001; E = 164.4;  First = $false },
    @{ Name = "British Columbia";        D = 155.6;  E = 130.5;              First = $false }
)

$row = 272
foreach ($p in $provinceRows) {
    $wsProvince.Cells.Item($row, 1).Value = $newDate
    $wsProvince.Cells.Item($row, 1).NumberFormat = $dateFmt

    $wsProvince.Cells.Item($row, 2).Value = $p.Name
    if ($p.First) {
        $wsProvince.Cells.Item($row, 2).NumberFormat = $dateFmt
    }

    $wsProvince.Cells.Item($row, 4).Value = $p.D
    $wsProvince.Cells.Item($row, 5).Value = $p.E

    $wsProvince.Range("C" + $row).Formula = "=(D" + $row + "-E" + $row + ")/E" + $row + "*100"

    $row = $row + 1
}

# ---------------------------------------------------------------------
# View state: update selections on both sheets. "Canada" is selected
# first so the final selection on "Province" leaves it as the active tab
# (matching the workbook's original active sheet).
# ---------------------------------------------------------------------
$wsCanada.Range("A29").Select()
$wsProvince.Range("D282").Select()

Write-Output "Added 2022-04 (published 9 May 2022) observations: 1 Canada row + 10 Province rows"
